# "Generate Report for Archive"
#
# The two handed-off documents have progressed out of "Ready for handoff"
# into "In Translation". Update the status cells on every sheet that
# surfaces that status (the Overview rollup plus each per-locale sheet),
# then re-fit the now-narrower status columns so the report reads cleanly.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-language status lives in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Re-fit the status columns now that "In Translation" is shorter than
#     "Ready for handoff". Feeding 12.5 here lands ColumnWidth on the same
#     rendered width the other (unmodified) narrow columns use. ---
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth

# --- quick sanity log ---
Write-Host "Overview!E2 -> $($wsOverview.Range('E2').Text)"
Write-Host "zh-cn!C2 -> $($wsZhCn.Range('C2').Text)"
Write-Host "de-de!C2 -> $($wsDeDe.Range('C2').Text)"
Write-Host "Overview col E width -> $($wsOverview.Columns.Item(5).ColumnWidth)"
